$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.532141
$ws.Range("H2").Value = 4.596423
$ws.Range("I2").Value = 0.08900664250669833
$ws.Range("J2").Value = 0.08900664250669831
$ws.Range("M2").Value = 1.294217333333333
$ws.Range("N2").Value = 3.882652
$ws.Range("O2").Value = 0.1864098899142058
$ws.Range("P2").Value = 0.1864098899142058
$ws.Range("Q2").Value = 1.982923439310667
$ws.Range("R2").Value = 17.846310953796
$ws.Range("S2").Value = 0.01659171843130671
$ws.Range("T2").Value = 0.0165917184313067
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.532141
$ws.Range("H3").Value = 4.596423
$ws.Range("I3").Value = 0.08900664250669833
$ws.Range("J3").Value = 0.08900664250669831
$ws.Range("O3").Value = 0.2342636243010983
$ws.Range("P3").Value = 0.2342636243010983
$ws.Range("Q3").Value = 2.491964518719
$ws.Range("R3").Value = 22.427680668471
$ws.Range("S3").Value = 0.02085101866049134
$ws.Range("T3").Value = 0.02085101866049133
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.532141
$ws.Range("H4").Value = 4.596423
$ws.Range("I4").Value = 0.08900664250669833
$ws.Range("J4").Value = 0.08900664250669831
$ws.Range("M4").Value = 3.488917666666667
$ws.Range("N4").Value = 10.466753
$ws.Range("O4").Value = 0.5025189675740148
$ws.Range("P4").Value = 0.5025189675740148
$ws.Range("Q4").Value = 5.345513802724334
$ws.Range("R4").Value = 48.109624224519
$ws.Range("S4").Value = 0.04472752609969547
$ws.Range("T4").Value = 0.04472752609969546
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.532141
$ws.Range("H5").Value = 4.596423
$ws.Range("I5").Value = 0.08900664250669833
$ws.Range("J5").Value = 0.08900664250669831
$ws.Range("M5").Value = 0.5332636666666667
$ws.Range("N5").Value = 1.599791
$ws.Range("O5").Value = 0.07680751821068107
$ws.Range("P5").Value = 0.07680751821068106
$ws.Range("Q5").Value = 0.8170351275103334
$ws.Range("R5").Value = 7.353316147593
$ws.Range("S5").Value = 0.006836379315204813
$ws.Range("T5").Value = 0.006836379315204809
$ws.Range("I6").Value = 0.6169137955113024
$ws.Range("J6").Value = 0.6169137955113023
$ws.Range("M6").Value = 1.294217333333333
$ws.Range("N6").Value = 3.882652
$ws.Range("O6").Value = 0.1864098899142058
$ws.Range("P6").Value = 0.1864098899142058
$ws.Range("Q6").Value = 13.74383743394667
$ws.Range("R6").Value = 123.69453690552
$ws.Range("S6").Value = 0.1149988327078168
$ws.Range("T6").Value = 0.1149988327078167
$ws.Range("I7").Value = 0.6169137955113024
$ws.Range("J7").Value = 0.6169137955113023
$ws.Range("O7").Value = 0.2342636243010983
$ws.Range("P7").Value = 0.2342636243010983
$ws.Range("S7").Value = 0.1445204616178243
$ws.Range("T7").Value = 0.1445204616178243
$ws.Range("I8").Value = 0.6169137955113024
$ws.Range("J8").Value = 0.6169137955113023
$ws.Range("M8").Value = 3.488917666666667
$ws.Range("N8").Value = 10.466753
$ws.Range("O8").Value = 0.5025189675740148
$ws.Range("P8").Value = 0.5025189675740148
$ws.Range("Q8").Value = 37.05028204775333
$ws.Range("R8").Value = 333.45253842978
$ws.Range("S8").Value = 0.3100108836025066
$ws.Range("T8").Value = 0.3100108836025066
$ws.Range("I9").Value = 0.6169137955113024
$ws.Range("J9").Value = 0.6169137955113023
$ws.Range("M9").Value = 0.5332636666666667
$ws.Range("N9").Value = 1.599791
$ws.Range("O9").Value = 0.07680751821068107
$ws.Range("P9").Value = 0.07680751821068106
$ws.Range("Q9").Value = 5.662950847073334
$ws.Range("R9").Value = 50.96655762366
$ws.Range("S9").Value = 0.04738361758315474
$ws.Range("T9").Value = 0.04738361758315472
$ws.Range("G10").Value = 4.902263666666666
$ws.Range("H10").Value = 14.706791
$ws.Range("I10").Value = 0.2847871244569372
$ws.Range("J10").Value = 0.2847871244569371
$ws.Range("M10").Value = 1.294217333333333
$ws.Range("N10").Value = 3.882652
$ws.Range("O10").Value = 0.1864098899142058
$ws.Range("P10").Value = 0.1864098899142058
$ws.Range("Q10").Value = 6.344594609970222
$ws.Range("R10").Value = 57.101351489732
$ws.Range("S10").Value = 0.0530871365190009
$ws.Range("T10").Value = 0.05308713651900088
$ws.Range("G11").Value = 4.902263666666666
$ws.Range("H11").Value = 14.706791
$ws.Range("I11").Value = 0.2847871244569372
$ws.Range("J11").Value = 0.2847871244569371
$ws.Range("O11").Value = 0.2342636243010983
$ws.Range("P11").Value = 0.2342636243010983
$ws.Range("Q11").Value = 7.973330861022998
$ws.Range("R11").Value = 71.759977749207
$ws.Range("S11").Value = 0.06671526392957004
$ws.Range("T11").Value = 0.06671526392957003
$ws.Range("G12").Value = 4.902263666666666
$ws.Range("H12").Value = 14.706791
$ws.Range("I12").Value = 0.2847871244569372
$ws.Range("J12").Value = 0.2847871244569371
$ws.Range("M12").Value = 3.488917666666667
$ws.Range("N12").Value = 10.466753
$ws.Range("O12").Value = 0.5025189675740148
$ws.Range("P12").Value = 0.5025189675740148
$ws.Range("Q12").Value = 17.10359431329145
$ws.Range("R12").Value = 153.932348819623
$ws.Range("S12").Value = 0.1431109317604725
$ws.Range("T12").Value = 0.1431109317604725
$ws.Range("G13").Value = 4.902263666666666
$ws.Range("H13").Value = 14.706791
$ws.Range("I13").Value = 0.2847871244569372
$ws.Range("J13").Value = 0.2847871244569371
$ws.Range("M13").Value = 0.5332636666666667
$ws.Range("N13").Value = 1.599791
$ws.Range("O13").Value = 0.07680751821068107
$ws.Range("P13").Value = 0.07680751821068106
$ws.Range("Q13").Value = 2.614199097853445
$ws.Range("R13").Value = 23.527791880681
$ws.Range("S13").Value = 0.0218737922478937
$ws.Range("T13").Value = 0.02187379224789369
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.159958
$ws.Range("H14").Value = 0.479874
$ws.Range("I14").Value = 0.009292437525062282
$ws.Range("J14").Value = 0.009292437525062281
$ws.Range("M14").Value = 1.294217333333333
$ws.Range("N14").Value = 3.882652
$ws.Range("O14").Value = 0.1864098899142058
$ws.Range("P14").Value = 0.1864098899142058
$ws.Range("Q14").Value = 0.2070204162053333
$ws.Range("R14").Value = 1.863183745848
$ws.Range("S14").Value = 0.001732202256081496
$ws.Range("T14").Value = 0.001732202256081495
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.159958
$ws.Range("H15").Value = 0.479874
$ws.Range("I15").Value = 0.009292437525062282
$ws.Range("J15").Value = 0.009292437525062281
$ws.Range("O15").Value = 0.2342636243010983
$ws.Range("P15").Value = 0.2342636243010983
$ws.Range("Q15").Value = 0.260165128722
$ws.Range("R15").Value = 2.341486158498
$ws.Range("S15").Value = 0.002176880093212618
$ws.Range("T15").Value = 0.002176880093212617
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.159958
$ws.Range("H16").Value = 0.479874
$ws.Range("I16").Value = 0.009292437525062282
$ws.Range("J16").Value = 0.009292437525062281
$ws.Range("M16").Value = 3.488917666666667
$ws.Range("N16").Value = 10.466753
$ws.Range("O16").Value = 0.5025189675740148
$ws.Range("P16").Value = 0.5025189675740148
$ws.Range("Q16").Value = 0.5580802921246667
$ws.Range("R16").Value = 5.022722629122
$ws.Range("S16").Value = 0.004669626111340332
$ws.Range("T16").Value = 0.004669626111340331
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.159958
$ws.Range("H17").Value = 0.479874
$ws.Range("I17").Value = 0.009292437525062282
$ws.Range("J17").Value = 0.009292437525062281
$ws.Range("M17").Value = 0.5332636666666667
$ws.Range("N17").Value = 1.599791
$ws.Range("O17").Value = 0.07680751821068107
$ws.Range("P17").Value = 0.07680751821068106
$ws.Range("Q17").Value = 0.08529978959266667
$ws.Range("R17").Value = 0.7676981063339999
$ws.Range("S17").Value = 0.0007137290644278375
$ws.Range("T17").Value = 0.0007137290644278372
